$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Sınıf Sayısı" (class count) values in column I: every data row (2-29)
# should be 2 instead of 1.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 9).Value = 2
}

# Give column I an explicit best-fit width, matching the other data columns.
$ws.Columns.Item(9).AutoFit() | Out-Null

# Reflect the updated selection left behind after the fix.
$ws.Range("I2:I29").Select()
